$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E1:E3").ClearContents() | Out-Null
$ws.Columns("E:E").Select() | Out-Null
